$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 184.5
$ws.Range("I101").Value = 181.4
$ws.Range("J101").Value = 200
$ws.Range("K101").Value = 544.2
$ws.Range("L101").Value = 600
$ws.Range("M101").Value = 1077.8
$ws.Range("N101").Value = -3844
$ws.Range("H132").Value = 4296.2104
$ws.Range("I132").Value = 4277.027
$ws.Range("K132").Value = 12831.081
$ws.Range("M132").Value = -10301.081
$ws.Range("H133").Value = 49911
$ws.Range("J133").Value = 49911
$ws.Range("L133").Value = 49911
$ws.Range("N133").Value = -60031
$ws.Range("H135").Value = 1552.7742
$ws.Range("I135").Value = 1560.1538
$ws.Range("K135").Value = 14041.3842
$ws.Range("M135").Value = -11506.3842
$ws.Range("H137").Value = 12477.685
$ws.Range("I137").Value = 1319
$ws.Range("J137").Value = 20593.092
$ws.Range("K137").Value = 3957
$ws.Range("L137").Value = 61779.276
$ws.Range("M137").Value = -1407
$ws.Range("N137").Value = -66879.276
$ws.Range("H138").Value = 5451.933
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 985.4935
$ws.Range("I32").Value = 984.9459000000001
$ws.Range("K32").Value = 984.9459000000001
$ws.Range("M32").Value = -697.9459000000001
$ws.Range("H61").Value = 4128.375
$ws.Range("I61").Value = 4128.375
$ws.Range("K61").Value = 4128.375
$ws.Range("M61").Value = -3916.375
$ws.Range("H74").Value = 292830.06
$ws.Range("I74").Value = 313719.5
$ws.Range("K74").Value = 313719.5
$ws.Range("M74").Value = -312845.5
$ws.Range("H77").Value = 292830.06
$ws.Range("I77").Value = 313719.5
$ws.Range("K77").Value = 1568597.5
$ws.Range("M77").Value = -1564229.5
$ws.Range("H97").Value = 1376.7273
$ws.Range("I97").Value = 1345.9656
$ws.Range("K97").Value = 1345.9656
$ws.Range("M97").Value = -849.9656
$ws.Range("H132").Value = 255117.48
$ws.Range("I132").Value = 279854.34
$ws.Range("K132").Value = 839563.02
$ws.Range("M132").Value = -837033.02
$ws.Range("H136").Value = 4128.375
$ws.Range("I136").Value = 4128.375
$ws.Range("K136").Value = 12385.125
$ws.Range("M136").Value = -9835.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = $null
$ws.Range("H86").Value = 588298.2
$ws.Range("I86").Value = 774560.7
$ws.Range("K86").Value = 774560.7
$ws.Range("M86").Value = -773437.7
$ws.Range("H89").Value = 588298.2
$ws.Range("I89").Value = 774560.7
$ws.Range("K89").Value = 3872803.5
$ws.Range("M89").Value = -3867187.5
$ws.Range("H134").Value = 104930.18
$ws.Range("I134").Value = 5406.3335
$ws.Range("J134").Value = 142251.62
$ws.Range("K134").Value = 16219.0005
$ws.Range("L134").Value = 426754.86
$ws.Range("M134").Value = -13684.0005
$ws.Range("N134").Value = -431824.86

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 757832.25
$ws.Range("I31").Value = 1413544.8
$ws.Range("K31").Value = 1413544.8
$ws.Range("M31").Value = -1413249.8
$ws.Range("H34").Value = 757832.25
$ws.Range("I34").Value = 1413544.8
$ws.Range("K34").Value = 1413544.8
$ws.Range("M34").Value = -1413342.8
$ws.Range("H50").Value = 23694
$ws.Range("J50").Value = 25507.428
$ws.Range("L50").Value = 25507.428
$ws.Range("N50").Value = -26757.428
$ws.Range("H51").Value = 18189.182
$ws.Range("J51").Value = 51099
$ws.Range("L51").Value = 51099
$ws.Range("N51").Value = -52571
$ws.Range("H59").Value = 30000
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = $null
$ws.Range("H60").Value = 35551
$ws.Range("J60").Value = 51103
$ws.Range("L60").Value = 51103
$ws.Range("N60").Value = -52125
$ws.Range("H61").Value = 18189.182
$ws.Range("J61").Value = 51099
$ws.Range("L61").Value = 51099
$ws.Range("N61").Value = -51795
$ws.Range("H68").Value = 79736.8
$ws.Range("J68").Value = 79600
$ws.Range("L68").Value = 79600
$ws.Range("N68").Value = -81098
$ws.Range("H71").Value = 79736.8
$ws.Range("J71").Value = 79600
$ws.Range("L71").Value = 238800
$ws.Range("N71").Value = -246288

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 497.83334
$ws.Range("I2").Value = 92.2
$ws.Range("J2").Value = 787.5714
$ws.Range("K2").Value = 553.2
$ws.Range("L2").Value = 4725.428400000001
$ws.Range("M2").Value = -440.2
$ws.Range("N2").Value = -4951.428400000001
$ws.Range("H4").Value = 14726713
$ws.Range("I4").Value = 824323.8
$ws.Range("K4").Value = 2472971.4
$ws.Range("M4").Value = -2472859.4
$ws.Range("H12").Value = 195.26666
$ws.Range("J12").Value = 217.61539
$ws.Range("L12").Value = 652.84617
$ws.Range("N12").Value = -998.84617
$ws.Range("H112").Value = 145555.14
$ws.Range("I112").Value = 145555.14
$ws.Range("K112").Value = 436665.42
$ws.Range("M112").Value = -435557.42
$ws.Range("H114").Value = 604.125
$ws.Range("I114").Value = 433.66666
$ws.Range("K114").Value = 1300.99998
$ws.Range("M114").Value = 1953.00002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 55995
$ws.Range("I62").Value = 55990
$ws.Range("K62").Value = 55990
$ws.Range("M62").Value = -55304
$ws.Range("H63").Value = 39701.332
$ws.Range("I63").Value = 30104
$ws.Range("J63").Value = 44500
$ws.Range("K63").Value = 30104
$ws.Range("L63").Value = 44500
$ws.Range("N63").Value = -45872
$ws.Range("M63").Value = -29418
$ws.Range("H65").Value = 55995
$ws.Range("I65").Value = 55990
$ws.Range("K65").Value = 167970
$ws.Range("M65").Value = -164538
$ws.Range("H66").Value = 39701.332
$ws.Range("I66").Value = 30104
$ws.Range("J66").Value = 44500
$ws.Range("K66").Value = 90312
$ws.Range("L66").Value = 133500
$ws.Range("N66").Value = -140364
$ws.Range("M66").Value = -86880
$ws.Range("H122").Value = 411737.62
$ws.Range("I122").Value = 528117.5600000001
$ws.Range("K122").Value = 1584352.68
$ws.Range("M122").Value = -1581902.68
$ws.Range("H132").Value = 51928.777
$ws.Range("I132").Value = 18115.428
$ws.Range("K132").Value = 54346.284
$ws.Range("M132").Value = -51816.284
$ws.Range("H135").Value = 91666.664
$ws.Range("J135").Value = 91666.664
$ws.Range("L135").Value = 91666.664
$ws.Range("N135").Value = -101806.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 9350
$ws.Range("I10").Value = 9350
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 9350
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -9210
$ws.Range("N10").Value = $null
$ws.Range("H63").Value = 48330
$ws.Range("J63").Value = 47500
$ws.Range("L63").Value = 47500
$ws.Range("N63").Value = -48998
$ws.Range("H66").Value = 48330
$ws.Range("J66").Value = 47500
$ws.Range("L66").Value = 142500
$ws.Range("N66").Value = -149988
$ws.Range("H132").Value = 4910.0835
$ws.Range("I132").Value = 4152.4
$ws.Range("J132").Value = 6632.091
$ws.Range("K132").Value = 12457.2
$ws.Range("L132").Value = 19896.273
$ws.Range("M132").Value = -9927.199999999999
$ws.Range("N132").Value = -24956.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 109500
$ws.Range("J46").Value = 109500
$ws.Range("L46").Value = 109500
$ws.Range("N46").Value = -109962
$ws.Range("H63").Value = 75000
$ws.Range("J63").Value = 75000
$ws.Range("L63").Value = 75000
$ws.Range("N63").Value = -76248
$ws.Range("H64").Value = 77998
$ws.Range("I64").Value = 74990
$ws.Range("K64").Value = 74990
$ws.Range("M64").Value = -74742
$ws.Range("H66").Value = 75000
$ws.Range("J66").Value = 75000
$ws.Range("L66").Value = 225000
$ws.Range("N66").Value = -231240
$ws.Range("H67").Value = 77998
$ws.Range("I67").Value = 74990
$ws.Range("K67").Value = 74990
$ws.Range("M67").Value = -74132
$ws.Range("H100").Value = 951
$ws.Range("I100").Value = 951
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1902
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1361
$ws.Range("N100").Value = $null
$ws.Range("H132").Value = 32873.082
$ws.Range("J132").Value = 75130.31
$ws.Range("L132").Value = 225390.93
$ws.Range("N132").Value = -230450.93
$ws.Range("H134").Value = 109500
$ws.Range("J134").Value = 109500
$ws.Range("L134").Value = 328500
$ws.Range("N134").Value = -333570
$ws.Range("H136").Value = 277555.62
$ws.Range("I136").Value = 265301.12
$ws.Range("J136").Value = 355167.5
$ws.Range("K136").Value = 795903.36
$ws.Range("L136").Value = 1065502.5
$ws.Range("M136").Value = -793353.36
$ws.Range("N136").Value = -1070602.5
